$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 189
$ws.Range("B189").Value = 7952778
$ws.Range("E189").Value = "Sloga"
$ws.Range("F189").Value = "Siroki Brijeg"
$ws.Range("G189").Value = 2
$ws.Range("H189").Value = 3
$ws.Range("J189").Value = 2
$ws.Range("K189").Value = "A"
$ws.Range("L189").Value = 1.727
$ws.Range("M189").Value = 3.75
$ws.Range("N189").Value = 3.75
$ws.Range("O189").Value = 1.7
$ws.Range("P189").Value = 3.9
$ws.Range("Q189").Value = 3.9
$ws.Range("R189").Value = -0.75
$ws.Range("U189").Value = 2.25
$ws.Range("V189").Value = 1.8
$ws.Range("W189").Value = 2
$ws.Range("X189").Value = -1
$ws.Range("Z189").Value = 2.9
$ws.Range("AA189").Value = -1
$ws.Range("AB189").Value = 0.825
$ws.Range("AC189").Value = 0.8
$ws.Range("AD189").Value = -1

# Row 190
$ws.Range("B190").Value = 7952781
$ws.Range("E190").Value = "Zvijezda 09"
$ws.Range("F190").Value = "Zeljeznicar"
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 5
$ws.Range("I190").Value = 0
$ws.Range("K190").Value = "A"
$ws.Range("L190").Value = 2.15
$ws.Range("M190").Value = 3.25
$ws.Range("N190").Value = 2.9
$ws.Range("O190").Value = 3.6
$ws.Range("P190").Value = 3.4
$ws.Range("Q190").Value = 1.85
$ws.Range("R190").Value = 0.5
$ws.Range("S190").Value = 1.875
$ws.Range("T190").Value = 1.925
$ws.Range("U190").Value = 2.5
$ws.Range("V190").Value = 1.975
$ws.Range("W190").Value = 1.825
$ws.Range("Y190").Value = -1
$ws.Range("Z190").Value = 0.8500000000000001
$ws.Range("AB190").Value = 0.925
$ws.Range("AC190").Value = 0.9750000000000001

# Row 191
$ws.Range("B191").Value = 7952776
$ws.Range("E191").Value = "FK Sarajevo"
$ws.Range("F191").Value = "NK Posusje"
$ws.Range("G191").Value = 1
$ws.Range("H191").Value = 1
$ws.Range("I191").Value = 0
$ws.Range("J191").Value = 0
$ws.Range("K191").Value = "D"
$ws.Range("L191").Value = 1.571
$ws.Range("M191").Value = 3.4
$ws.Range("N191").Value = 5.5
$ws.Range("O191").Value = 1.363
$ws.Range("P191").Value = 3.9
$ws.Range("Q191").Value = 8
$ws.Range("R191").Value = -1.25
$ws.Range("S191").Value = 1.85
$ws.Range("T191").Value = 1.95
$ws.Range("U191").Value = 2.75
$ws.Range("V191").Value = 1.925
$ws.Range("W191").Value = 1.875
$ws.Range("X191").Value = -1
$ws.Range("Y191").Value = 2.9
$ws.Range("AB191").Value = 0.95
$ws.Range("AC191").Value = -1
$ws.Range("AD191").Value = 0.875

# Row 192
$ws.Range("B192").Value = 7952777
$ws.Range("E192").Value = "Borac Banja Luka"
$ws.Range("F192").Value = "NK Igman Konjic"
$ws.Range("G192").Value = 4
$ws.Range("I192").Value = 1
$ws.Range("K192").Value = "H"
$ws.Range("L192").Value = 1.25
$ws.Range("M192").Value = 5.75
$ws.Range("N192").Value = 7
$ws.Range("O192").Value = 1.2
$ws.Range("P192").Value = 5.75
$ws.Range("Q192").Value = 12
$ws.Range("R192").Value = -2
$ws.Range("S192").Value = 1.95
$ws.Range("T192").Value = 1.85
$ws.Range("U192").Value = 3.25
$ws.Range("V192").Value = 1.9
$ws.Range("W192").Value = 1.9
$ws.Range("X192").Value = 0.2
$ws.Range("Z192").Value = -1
$ws.Range("AB192").Value = 0.8500000000000001
$ws.Range("AC192").Value = 0.8999999999999999

# Row 193
$ws.Range("B193").Value = 7952780
$ws.Range("E193").Value = "Velez Mostar"
$ws.Range("F193").Value = "GOSK Gabela"
$ws.Range("G193").Value = 3
$ws.Range("H193").Value = 3
$ws.Range("I193").Value = 1
$ws.Range("K193").Value = "D"
$ws.Range("L193").Value = 1.4
$ws.Range("M193").Value = 4
$ws.Range("N193").Value = 7
$ws.Range("O193").Value = 1.363
$ws.Range("P193").Value = 4.2
$ws.Range("Q193").Value = 8
$ws.Range("R193").Value = -1.5
$ws.Range("S193").Value = 2
$ws.Range("T193").Value = 1.8
$ws.Range("U193").Value = 2.75
$ws.Range("V193").Value = 1.825
$ws.Range("W193").Value = 1.975
$ws.Range("Y193").Value = 3.2
$ws.Range("Z193").Value = -1
$ws.Range("AB193").Value = 0.8
$ws.Range("AC193").Value = 0.825

# Row 194
$ws.Range("B194").Value = 7952779
$ws.Range("E194").Value = "Zrinjski Mostar"
$ws.Range("F194").Value = "FK Tuzla City"
$ws.Range("G194").Value = 4
$ws.Range("H194").Value = 0
$ws.Range("I194").Value = 2
$ws.Range("K194").Value = "H"
$ws.Range("L194").Value = 1.25
$ws.Range("M194").Value = 5.75
$ws.Range("N194").Value = 7
$ws.Range("O194").Value = 1.055
$ws.Range("P194").Value = 13
$ws.Range("Q194").Value = 17
$ws.Range("R194").Value = -3.5
$ws.Range("S194").Value = 1.975
$ws.Range("T194").Value = 1.825
$ws.Range("U194").Value = 4.75
$ws.Range("V194").Value = 1.825
$ws.Range("W194").Value = 1.975
$ws.Range("X194").Value = 0.05499999999999994
$ws.Range("Y194").Value = -1
$ws.Range("AA194").Value = 0.9750000000000001
$ws.Range("AB194").Value = -1
